$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Overall" row (currently row 21), shifting the
# Overall/spacer/footer rows down by one.
$ws.Rows.Item(21).Insert()

# Copy formatting from the row above (row 20, the last weighted KPI row) so
# the new row matches the existing KPI rows' look (fonts/fills/alignment).
$ws.Range("A20:G20").Copy()
$ws.Range("A21:G21").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(21).RowHeight = $ws.Rows.Item(20).RowHeight

# Populate the new KPI row with "Client love" data.
$ws.Cells.Item(21, 1).Value = "Client love"
$ws.Cells.Item(21, 3).Value = 0.2
$ws.Cells.Item(21, 5).Formula = "=PRODUCT(C21:D21) * 10"

# Fix up the "Overall" row's formulas (now row 22) so the sum includes the
# newly inserted row.
$ws.Cells.Item(22, 5).Formula = "=SUM(E3:E21)"
$ws.Cells.Item(22, 7).Formula = "= E22/F22"

# Match the author's final selection (as seen in the diff). (Note: this
# runtime does not persist a separate window scroll/topLeftCell position,
# only the active selection, so that's all we can reproduce here.)
[void]$ws.Range("A17").Select()

$wb.Save()
